$p = $ppt.ActivePresentation

# --- 1) Footer "last modified" date fields: 28/01/2021 -> 29/01/2021 ---
# Slide master + all of its custom layouts hold a literal "Date Placeholder"
# shape whose TextRange we can rewrite directly.
$master = $p.SlideMaster

$masterShapes = $master.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $shp = $masterShapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "29/01/2021"
    }
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $lay = $layouts.Item($L)
    $layShapes = $lay.Shapes
    for ($i = 1; $i -le $layShapes.Count; $i++) {
        $shp = $layShapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "29/01/2021"
        }
    }
}

# The Notes Master's own shape collection aliases back onto the slide
# master in this host, so update it through HeadersFooters instead.
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "29/01/2021"

# --- 2) Slide 18 ("Components"): emphasise a handful of servlet / page
#        bullets that are now wired up, per the commit message. ---
$slide = $p.Slides.Item(18)

# "Client components" bullet list (servlets)
$servletShape = $slide.Shapes.Item(2)
$servletText = $servletShape.TextFrame.TextRange

$servletText.Paragraphs(4, 1).Font.Bold = $true
$servletText.Paragraphs(5, 1).Font.Bold = $true
$servletText.Paragraphs(6, 1).Font.Underline = $true
$servletText.Paragraphs(7, 1).Font.Underline = $true
$servletText.Paragraphs(8, 1).Font.Underline = $true

# "Views" bullet list (html pages)
$viewsShape = $slide.Shapes.Item(3)
$viewsText = $viewsShape.TextFrame.TextRange

$viewsText.Paragraphs(4, 1).Font.Bold = $true
$viewsText.Paragraphs(5, 1).Font.Bold = $true

$greetingsPara = $viewsText.Paragraphs(6, 1)
$greetingsPara.Characters(1, 9).Font.Underline = $true

Write-Output "edits applied"
